$d = $word.ActiveDocument

# 1. Name casing change
$d.Content.Find.Execute("DHEERAJ CHAND", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dheeraj Chand", 2)

# 2. Title placeholder
$d.Content.Find.Execute("Senior Product Marketing Manager", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Professional Title", 2)

# 3. Contact line formatting
$d.Content.Find.Execute("(202) 550-7110 | Dheeraj.Chand@gmail.com", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "202.550.7110 | dheeraj.chand@gmail.com", 2)

# 4. Years of experience wording
$d.Content.Find.Execute("20+ years of experience", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "21 years of experience", 2)

# 5. Company name -> placeholder
$d.Content.Find.Execute("Siege Analytics, Austin, TX | 2005", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Your Company Name, Your City, ST | 2005", 2)

# 6. Remove product names from bullet
$d.Content.Find.Execute("multiple SaaS platform launches including BALLISTA and DAMON, achieving", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "multiple SaaS platform launches, achieving", 2)

# 7. Remove product names from achievements bullet
$d.Content.Find.Execute("multiple B2B SaaS platforms (BALLISTA, DAMON, SimCrisis, RACSO) used by", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "multiple B2B SaaS platforms used by", 2)

# 8. Delete the trailing "Market Intelligence & Research Leadership" / "Cross-Functional
#    Leadership & Collaboration" sections (last block of the document).
$startP = $null
$endP = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Market Intelligence*Research Leadership*") {
        $startP = $i
    }
}
if ($startP -ne $null) {
    $rStart = $d.Paragraphs.Item($startP).Range.Start
    $rEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
    $r = $d.Range($rStart, $rEnd)
    $r.Delete()
}

# 9. Delete the multi-job history block: from "DATA PRODUCTS MANAGER" through the
#    "RESEARCH DIRECTOR & PRODUCT MANAGER" entry's last bullet, right before
#    "KEY ACHIEVEMENTS AND IMPACT".
$startP = $null
$endP = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "DATA PRODUCTS MANAGER*") {
        $startP = $i
    }
    if ($t -like "KEY ACHIEVEMENTS AND IMPACT*") {
        $endP = $i
    }
}
if (($startP -ne $null) -and ($endP -ne $null)) {
    $rStart = $d.Paragraphs.Item($startP).Range.Start
    $rEnd = $d.Paragraphs.Item($endP).Range.Start
    $r = $d.Range($rStart, $rEnd)
    $r.Delete()
}
